$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "26.263.25"
$ws.Cells.Item(2, 5).Value = "  -0.37%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.595.31"
$ws.Cells.Item(3, 5).Value = "  +0.20%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "212.85"
$ws.Cells.Item(5, 5).Value = "  +0.70%  "

$ws.Cells.Item(6, 5).Value = "  -0.73%  "

$ws.Cells.Item(7, 5).Value = "  +0.00%  "

$ws.Cells.Item(8, 5).Value = "  -0.32%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.0607"
$ws.Cells.Item(9, 5).Value = "  -0.56%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "18.97"
$ws.Cells.Item(10, 5).Value = "  -2.55%  "

$ws.Cells.Item(11, 5).Value = "  +0.51%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.819.61"
$ws.Cells.Item(12, 5).Value = "  +0.19%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.591.10"
$ws.Cells.Item(13, 5).Value = "  -0.14%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.00"
$ws.Cells.Item(14, 5).Value = "  -1.31%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.509"
$ws.Cells.Item(15, 5).Value = "  -2.62%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "63.96"
$ws.Cells.Item(16, 5).Value = "  -1.23%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "26.256.37"
$ws.Cells.Item(17, 5).Value = "  -0.39%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.0₃0726"
$ws.Cells.Item(18, 5).Value = "  -0.63%  "

$ws.Cells.Item(19, 2).Value = "BitcoinCash"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "214.08"
$ws.Cells.Item(19, 5).Value = "  +0.97%  "

$ws.Cells.Item(20, 2).Value = "Chainlink"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.35"
$ws.Cells.Item(20, 5).Value = "  -1.88%  "

$ws.Cells.Item(21, 5).Value = "  -0.08%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.30"
$ws.Cells.Item(22, 5).Value = "  -0.14%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "9.06"
$ws.Cells.Item(23, 5).Value = "  +0.49%  "

$ws.Cells.Item(24, 5).Value = "  -3.09%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "145.04"
$ws.Cells.Item(25, 5).Value = "  +0.14%  "

$ws.Cells.Item(26, 5).Value = "  -0.01%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "6.97"
$ws.Cells.Item(27, 5).Value = "  -1.63%  "

$ws.Cells.Item(28, 5).Value = "  -0.62%  "

$ws.Cells.Item(29, 5).Value = "  -0.60%  "

$ws.Cells.Item(30, 5).Value = "  -2.58%  "

$ws.Cells.Item(31, 5).Value = "  +0.33%  "

$ws.Cells.Item(32, 5).Value = "  -0.54%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.420.66"
$ws.Cells.Item(33, 5).Value = "  +6.10%  "

$ws.Cells.Item(34, 5).Value = "  -0.16%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.43"
$ws.Cells.Item(35, 5).Value = "  -0.57%  "

$ws.Cells.Item(36, 5).Value = "  -1.00%  "

$ws.Cells.Item(37, 5).Value = "  -2.88%  "

$ws.Cells.Item(38, 5).Value = "  -0.84%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.822"
$ws.Cells.Item(39, 5).Value = "  +0.63%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "5.81"
$ws.Cells.Item(40, 5).Value = "  +0.68%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.969"
$ws.Cells.Item(42, 5).Value = "  -8.86%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.16"
$ws.Cells.Item(43, 5).Value = "  +0.97%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.765"
$ws.Cells.Item(44, 5).Value = "  +0.01%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.731.05"
$ws.Cells.Item(45, 5).Value = "  +0.11%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "60.98"
$ws.Cells.Item(46, 5).Value = "  -1.06%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "86.98"
$ws.Cells.Item(47, 5).Value = "  -1.22%  "

$ws.Cells.Item(48, 5).Value = "  -0.44%  "

$ws.Cells.Item(49, 5).Value = "  -0.65%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0958"
$ws.Cells.Item(50, 5).Value = "  -2.68%  "
